$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1000
$ws.Cells.Item(32, 10).Value = 1000
$ws.Cells.Item(32, 12).Value = 1000
$ws.Cells.Item(32, 14).Value = -1652
$ws.Cells.Item(74, 8).Value = 6327.1665
$ws.Cells.Item(74, 9).Value = 6327.1665
$ws.Cells.Item(74, 11).Value = 6327.1665
$ws.Cells.Item(74, 13).Value = -5391.1665
$ws.Cells.Item(77, 8).Value = 6327.1665
$ws.Cells.Item(77, 9).Value = 6327.1665
$ws.Cells.Item(77, 11).Value = 31635.8325
$ws.Cells.Item(77, 13).Value = -26955.8325
$ws.Cells.Item(86, 8).Value = 8399
$ws.Cells.Item(86, 9).Value = 8373.75
$ws.Cells.Item(86, 10).Value = 8500
$ws.Cells.Item(86, 11).Value = 8373.75
$ws.Cells.Item(86, 12).Value = 8500
$ws.Cells.Item(86, 13).Value = -7250.75
$ws.Cells.Item(86, 14).Value = -10746
$ws.Cells.Item(89, 8).Value = 8399
$ws.Cells.Item(89, 9).Value = 8373.75
$ws.Cells.Item(89, 10).Value = 8500
$ws.Cells.Item(89, 11).Value = 41868.75
$ws.Cells.Item(89, 12).Value = 42500
$ws.Cells.Item(89, 13).Value = -36252.75
$ws.Cells.Item(89, 14).Value = -53732
$ws.Cells.Item(121, 8).Value = 3000
$ws.Cells.Item(121, 10).Value = 3000
$ws.Cells.Item(121, 12).Value = 9000
$ws.Cells.Item(121, 14).Value = -12494
$ws.Cells.Item(138, 8).Value = 9177.799999999999
$ws.Cells.Item(138, 10).Value = 9392.579
$ws.Cells.Item(138, 12).Value = 28177.737
$ws.Cells.Item(138, 14).Value = -38457.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1330.7142
$ws.Cells.Item(2, 9).Value = 1331.909
$ws.Cells.Item(2, 10).Value = 1326.3334
$ws.Cells.Item(2, 11).Value = 1331.909
$ws.Cells.Item(2, 12).Value = 1326.3334
$ws.Cells.Item(2, 13).Value = -1218.909
$ws.Cells.Item(2, 14).Value = -1552.3334
$ws.Cells.Item(61, 8).Value = 2501.1667
$ws.Cells.Item(61, 9).Value = 2501.1667
$ws.Cells.Item(61, 11).Value = 2501.1667
$ws.Cells.Item(61, 13).Value = -2289.1667
$ws.Cells.Item(63, 8).Value = 15427.143
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 15427.143
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 15427.143
$ws.Cells.Item(63, 13).ClearContents()
$ws.Cells.Item(63, 14).Value = -16799.143
$ws.Cells.Item(66, 8).Value = 15427.143
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 15427.143
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 77135.715
$ws.Cells.Item(66, 13).ClearContents()
$ws.Cells.Item(66, 14).Value = -83999.715
$ws.Cells.Item(74, 8).Value = 5966.778
$ws.Cells.Item(74, 9).Value = 2962.625
$ws.Cells.Item(74, 11).Value = 2962.625
$ws.Cells.Item(74, 13).Value = -2088.625
$ws.Cells.Item(77, 8).Value = 5966.778
$ws.Cells.Item(77, 9).Value = 2962.625
$ws.Cells.Item(77, 11).Value = 14813.125
$ws.Cells.Item(77, 13).Value = -10445.125
$ws.Cells.Item(116, 8).Value = 1330.7142
$ws.Cells.Item(116, 9).Value = 1331.909
$ws.Cells.Item(116, 10).Value = 1326.3334
$ws.Cells.Item(116, 11).Value = 1331.909
$ws.Cells.Item(116, 12).Value = 1326.3334
$ws.Cells.Item(116, 13).Value = 962.0909999999999
$ws.Cells.Item(116, 14).Value = -5914.3334
$ws.Cells.Item(136, 8).Value = 2501.1667
$ws.Cells.Item(136, 9).Value = 2501.1667
$ws.Cells.Item(136, 11).Value = 7503.500100000001
$ws.Cells.Item(136, 13).Value = -4953.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1330.7142
$ws.Cells.Item(3, 9).Value = 1331.909
$ws.Cells.Item(3, 10).Value = 1326.3334
$ws.Cells.Item(3, 11).Value = 1331.909
$ws.Cells.Item(3, 12).Value = 1326.3334
$ws.Cells.Item(3, 13).Value = -1217.909
$ws.Cells.Item(3, 14).Value = -1554.3334
$ws.Cells.Item(86, 8).Value = 6059.8
$ws.Cells.Item(86, 9).Value = 4933.3335
$ws.Cells.Item(86, 10).Value = 6542.5713
$ws.Cells.Item(86, 11).Value = 4933.3335
$ws.Cells.Item(86, 12).Value = 6542.5713
$ws.Cells.Item(86, 13).Value = -3810.3335
$ws.Cells.Item(86, 14).Value = -8788.5713
$ws.Cells.Item(89, 8).Value = 6059.8
$ws.Cells.Item(89, 9).Value = 4933.3335
$ws.Cells.Item(89, 10).Value = 6542.5713
$ws.Cells.Item(89, 11).Value = 24666.6675
$ws.Cells.Item(89, 12).Value = 32712.8565
$ws.Cells.Item(89, 13).Value = -19050.6675
$ws.Cells.Item(89, 14).Value = -43944.85649999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 4000721
$ws.Cells.Item(8, 9).Value = 4000721
$ws.Cells.Item(8, 11).Value = 12002163
$ws.Cells.Item(8, 13).Value = -12002024
$ws.Cells.Item(119, 8).Value = 1750
$ws.Cells.Item(119, 9).Value = 1750
$ws.Cells.Item(119, 11).Value = 5250
$ws.Cells.Item(119, 13).Value = -412
$ws.Cells.Item(121, 8).Value = 674.75
$ws.Cells.Item(121, 9).Value = 649
$ws.Cells.Item(121, 10).Value = 683.3333
$ws.Cells.Item(121, 11).Value = 1947
$ws.Cells.Item(121, 12).Value = 2049.9999
$ws.Cells.Item(121, 13).Value = -637
$ws.Cells.Item(121, 14).Value = -4669.9999
$ws.Cells.Item(126, 8).Value = 4690
$ws.Cells.Item(126, 10).Value = 6600
$ws.Cells.Item(126, 12).Value = 19800
$ws.Cells.Item(126, 14).Value = -29680
$ws.Cells.Item(131, 8).Value = 1004.6
$ws.Cells.Item(131, 9).Value = 678
$ws.Cells.Item(131, 11).Value = 2034
$ws.Cells.Item(131, 13).Value = 3006
$ws.Cells.Item(137, 8).Value = 8500
$ws.Cells.Item(137, 10).Value = 2000
$ws.Cells.Item(137, 12).Value = 6000
$ws.Cells.Item(137, 14).Value = -16200
$ws.Cells.Item(139, 8).Value = 168287.5
$ws.Cells.Item(139, 9).Value = 168287.5
$ws.Cells.Item(139, 11).Value = 504862.5
$ws.Cells.Item(139, 13).Value = -499722.5
$ws.Cells.Item(140, 8).Value = 3079.6667
$ws.Cells.Item(140, 9).Value = 2664.4285
$ws.Cells.Item(140, 11).Value = 7993.2855
$ws.Cells.Item(140, 13).Value = -2813.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9000
$ws.Cells.Item(70, 9).Value = 9000
$ws.Cells.Item(70, 11).Value = 9000
$ws.Cells.Item(70, 13).Value = -8730
$ws.Cells.Item(73, 8).Value = 9000
$ws.Cells.Item(73, 9).Value = 9000
$ws.Cells.Item(73, 11).Value = 9000
$ws.Cells.Item(73, 13).Value = -8064
$ws.Cells.Item(80, 8).Value = 15167.5
$ws.Cells.Item(80, 9).Value = 12001.667
$ws.Cells.Item(80, 10).Value = 18333.334
$ws.Cells.Item(80, 11).Value = 12001.667
$ws.Cells.Item(80, 12).Value = 18333.334
$ws.Cells.Item(80, 13).Value = -11003.667
$ws.Cells.Item(80, 14).Value = -20329.334
$ws.Cells.Item(83, 8).Value = 15167.5
$ws.Cells.Item(83, 9).Value = 12001.667
$ws.Cells.Item(83, 10).Value = 18333.334
$ws.Cells.Item(83, 11).Value = 60008.335
$ws.Cells.Item(83, 12).Value = 91666.67
$ws.Cells.Item(83, 13).Value = -55016.335
$ws.Cells.Item(83, 14).Value = -101650.67
$ws.Cells.Item(122, 8).Value = 3606.2307
$ws.Cells.Item(122, 9).Value = 1988.1
$ws.Cells.Item(122, 10).Value = 9000
$ws.Cells.Item(122, 11).Value = 5964.299999999999
$ws.Cells.Item(122, 12).Value = 27000
$ws.Cells.Item(122, 13).Value = -3514.299999999999
$ws.Cells.Item(122, 14).Value = -31900
$ws.Cells.Item(132, 8).Value = 3006.9
$ws.Cells.Item(132, 9).Value = 1795.8572
$ws.Cells.Item(132, 11).Value = 5387.571599999999
$ws.Cells.Item(132, 13).Value = -2857.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2770.3333
$ws.Cells.Item(82, 9).Value = 733.2857
$ws.Cells.Item(82, 10).Value = 9900
$ws.Cells.Item(82, 11).Value = 733.2857
$ws.Cells.Item(82, 12).Value = 9900
$ws.Cells.Item(82, 13).Value = -372.2857
$ws.Cells.Item(82, 14).Value = -10622
$ws.Cells.Item(85, 8).Value = 2770.3333
$ws.Cells.Item(85, 9).Value = 733.2857
$ws.Cells.Item(85, 10).Value = 9900
$ws.Cells.Item(85, 11).Value = 733.2857
$ws.Cells.Item(85, 12).Value = 9900
$ws.Cells.Item(85, 13).Value = 514.7143
$ws.Cells.Item(85, 14).Value = -12396
$ws.Cells.Item(132, 8).Value = 2240.5386
$ws.Cells.Item(132, 9).Value = 1557.1818
$ws.Cells.Item(132, 10).Value = 5999
$ws.Cells.Item(132, 11).Value = 4671.5454
$ws.Cells.Item(132, 12).Value = 17997
$ws.Cells.Item(132, 13).Value = -2141.5454
$ws.Cells.Item(132, 14).Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 71001
$ws.Cells.Item(46, 10).Value = 71001
$ws.Cells.Item(46, 12).Value = 71001
$ws.Cells.Item(46, 14).Value = -71463
$ws.Cells.Item(134, 8).Value = 71001
$ws.Cells.Item(134, 10).Value = 71001
$ws.Cells.Item(134, 12).Value = 213003
$ws.Cells.Item(134, 14).Value = -218073
$ws.Cells.Item(136, 8).Value = 4650.8
$ws.Cells.Item(136, 9).Value = 4650.8
$ws.Cells.Item(136, 11).Value = 13952.4
$ws.Cells.Item(136, 13).Value = -11402.4
